$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 11) mirroring the existing rows' layout.
$ws.Range("A11").Value = 7
$ws.Range("B11").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C11").Value = "Ñuble"
$ws.Range("D11").Value = 44931
$ws.Range("D11").NumberFormat = $ws.Range("D10").NumberFormat
$ws.Range("E11").Value = 16
$ws.Range("F11").Value = "Fruta"
$ws.Range("G11").Value = 100101
$ws.Range("H11").Value = "Berries"
$ws.Range("I11").Value = 100101001
$ws.Range("J11").Value = "Arándano (blue)"
$ws.Range("K11").Value = "Sin especificar"
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 100
$ws.Range("N11").Value = 3000
$ws.Range("O11").Value = 3000
$ws.Range("P11").Value = 3000
$ws.Range("Q11").Value = "$/bandeja 2 kilos"
$ws.Range("R11").Value = "Provincia de Diguillín"
$ws.Range("S11").Value = 1500
$ws.Range("T11").Value = 2
